$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the BI/PD/UM/MP symbol columns (AE:AH) for rows 2,3,5 and row 6 (row4 unchanged: still "C")
$ws.Range("AE2:AH2").Value = "A"
$ws.Range("AE3:AH3").Value = "B"
$ws.Range("AE5:AH5").Value = "D"
$ws.Range("AE6:AH6").Value = "E"

# Update the selected cell in the sheet view
$ws.Range("B2").Select()
